$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10656.643
$ws.Range("I40").Value = 12779.4
$ws.Range("J40").Value = 5349.75
$ws.Range("K40").Value = 12779.4
$ws.Range("L40").Value = 5349.75
$ws.Range("M40").Value = -12604.4
$ws.Range("N40").Value = -5699.75
$ws.Range("H88").Value = 22276108.0
$ws.Range("J88").Value = 79768.9
$ws.Range("L88").Value = 79768.9
$ws.Range("N88").Value = -80580.9
$ws.Range("H91").Value = 22276108.0
$ws.Range("J91").Value = 79768.9
$ws.Range("L91").Value = 79768.9
$ws.Range("N91").Value = -82576.9
$ws.Range("H96").Value = 921.4286
$ws.Range("I96").Value = 741.6667
$ws.Range("K96").Value = 2225.0001
$ws.Range("M96").Value = -852.0001000000002
$ws.Range("H111").Value = 8337008.0
$ws.Range("I111").Value = 9618768.0
$ws.Range("K111").Value = 28856304.0
$ws.Range("M111").Value = -28853237.0
$ws.Range("H112").Value = 13066.375
$ws.Range("J112").Value = 13066.375
$ws.Range("L112").Value = 39199.125
$ws.Range("N112").Value = -41415.125
$ws.Range("H132").Value = 1698.409
$ws.Range("I132").Value = 1761.4147
$ws.Range("J132").Value = 837.3333
$ws.Range("K132").Value = 5284.2441
$ws.Range("L132").Value = 2511.9999
$ws.Range("M132").Value = -2754.2441
$ws.Range("N132").Value = -7571.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 5219.0
$ws.Range("J14").Value = 4828.5
$ws.Range("L14").Value = 4828.5
$ws.Range("N14").Value = -5178.5
$ws.Range("H32").Value = 1265008.9
$ws.Range("I32").Value = 1439180.6
$ws.Range("K32").Value = 1439180.6
$ws.Range("M32").Value = -1438893.6
$ws.Range("H45").Value = 2674.1
$ws.Range("I45").Value = 2506.7144
$ws.Range("J45").Value = 3064.6667
$ws.Range("K45").Value = 2506.7144
$ws.Range("L45").Value = 3064.6667
$ws.Range("M45").Value = -2129.7144
$ws.Range("N45").Value = -3818.6667
$ws.Range("H74").Value = 45332.344
$ws.Range("I74").Value = 58736.75
$ws.Range("K74").Value = 58736.75
$ws.Range("M74").Value = -57862.75
$ws.Range("H77").Value = 45332.344
$ws.Range("I77").Value = 58736.75
$ws.Range("K77").Value = 293683.75
$ws.Range("M77").Value = -289315.75
$ws.Range("H97").Value = 3211645.8
$ws.Range("I97").Value = 606.8333
$ws.Range("K97").Value = 606.8333
$ws.Range("M97").Value = -110.8333
$ws.Range("H132").Value = 1290347.2
$ws.Range("I132").Value = 2636953.2
$ws.Range("J132").Value = 7865.2856
$ws.Range("K132").Value = 7910859.600000001
$ws.Range("L132").Value = 23595.8568
$ws.Range("M132").Value = -7908329.600000001
$ws.Range("N132").Value = -28655.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 48910.0
$ws.Range("I50").Value = 0.0
$ws.Range("K50").Value = 0.0
$ws.Range("M50").ClearContents()
$ws.Range("H94").Value = 21278820.0
$ws.Range("I94").Value = 30304080.0
$ws.Range("K94").Value = 30304080.0
$ws.Range("M94").Value = -30303629.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7428.25
$ws.Range("I31").Value = 3111.95
$ws.Range("K31").Value = 3111.95
$ws.Range("M31").Value = -2816.95
$ws.Range("H34").Value = 7428.25
$ws.Range("I34").Value = 3111.95
$ws.Range("K34").Value = 3111.95
$ws.Range("M34").Value = -2909.95
$ws.Range("H43").Value = 33478.5
$ws.Range("J43").Value = 33478.5
$ws.Range("L43").Value = 33478.5
$ws.Range("N43").Value = -33846.5
$ws.Range("H58").Value = 6548.543
$ws.Range("I58").Value = 2097.647
$ws.Range("J58").Value = 10752.167
$ws.Range("K58").Value = 2097.647
$ws.Range("L58").Value = 10752.167
$ws.Range("M58").Value = -1894.647
$ws.Range("N58").Value = -11158.167
$ws.Range("H101").Value = 33478.5
$ws.Range("J101").Value = 33478.5
$ws.Range("L101").Value = 33478.5
$ws.Range("N101").Value = -39968.5
$ws.Range("H132").Value = 8198.0
$ws.Range("I132").Value = 3076.4
$ws.Range("K132").Value = 9229.2
$ws.Range("M132").Value = -6699.200000000001
$ws.Range("H134").Value = 7656.3184
$ws.Range("I134").Value = 3268.0908
$ws.Range("J134").Value = 12044.546
$ws.Range("K134").Value = 9804.2724
$ws.Range("L134").Value = 36133.638
$ws.Range("M134").Value = -7269.2724
$ws.Range("N134").Value = -41203.638
$ws.Range("H136").Value = 6548.543
$ws.Range("I136").Value = 2097.647
$ws.Range("J136").Value = 10752.167
$ws.Range("K136").Value = 6292.941
$ws.Range("L136").Value = 32256.501
$ws.Range("M136").Value = -3742.941
$ws.Range("N136").Value = -37356.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 217.9
$ws.Range("I23").Value = 194.66667
$ws.Range("J23").Value = 227.85715
$ws.Range("K23").Value = 584.00001
$ws.Range("L23").Value = 683.5714499999999
$ws.Range("M23").Value = -349.00001
$ws.Range("N23").Value = -1153.57145
$ws.Range("H68").Value = 4758.4443
$ws.Range("I68").Value = 2666.6667
$ws.Range("J68").Value = 5804.3335
$ws.Range("K68").Value = 8000.000100000001
$ws.Range("L68").Value = 17413.0005
$ws.Range("M68").Value = -7189.000100000001
$ws.Range("N68").Value = -19035.0005
$ws.Range("H71").Value = 4758.4443
$ws.Range("I71").Value = 2666.6667
$ws.Range("J71").Value = 5804.3335
$ws.Range("K71").Value = 24000.0003
$ws.Range("L71").Value = 52239.0015
$ws.Range("M71").Value = -19944.0003
$ws.Range("N71").Value = -60351.0015
$ws.Range("H80").Value = 29416334.0
$ws.Range("I80").Value = 20837374.0
$ws.Range("J80").Value = 50005840.0
$ws.Range("K80").Value = 62512122.0
$ws.Range("L80").Value = 150017520.0
$ws.Range("M80").Value = -62511186.0
$ws.Range("N80").Value = -150019392.0
$ws.Range("H83").Value = 29416334.0
$ws.Range("I83").Value = 20837374.0
$ws.Range("J83").Value = 50005840.0
$ws.Range("K83").Value = 187536366.0
$ws.Range("L83").Value = 450052560.0
$ws.Range("M83").Value = -187531686.0
$ws.Range("N83").Value = -450061920.0
$ws.Range("H87").Value = 41675960.0
$ws.Range("I87").Value = 142859860.0
$ws.Range("K87").Value = 428579580.0
$ws.Range("M87").Value = -428578332.0
$ws.Range("H90").Value = 41675960.0
$ws.Range("I90").Value = 142859860.0
$ws.Range("K90").Value = 1285738740.0
$ws.Range("M90").Value = -1285732500.0
$ws.Range("H114").Value = 457273.72
$ws.Range("I114").Value = 113.0
$ws.Range("J114").Value = 718508.44
$ws.Range("K114").Value = 339.0
$ws.Range("L114").Value = 2155525.32
$ws.Range("M114").Value = 2915.0
$ws.Range("N114").Value = -2162033.32
$ws.Range("H132").Value = 10451.053
$ws.Range("J132").Value = 12647.77
$ws.Range("L132").Value = 113829.93
$ws.Range("N132").Value = -118889.93
$ws.Range("H137").Value = 100585.71
$ws.Range("J137").Value = 136237.38
$ws.Range("L137").Value = 408712.14
$ws.Range("N137").Value = -418912.14
$ws.Range("H140").Value = 149945.38
$ws.Range("I140").Value = 183164.81
$ws.Range("K140").Value = 549494.4299999999
$ws.Range("M140").Value = -544314.4299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 998.0
$ws.Range("I43").Value = 998.0
$ws.Range("K43").Value = 998.0
$ws.Range("M43").Value = -847.0
$ws.Range("H102").Value = 1561.2
$ws.Range("I102").Value = 1311.7222
$ws.Range("K102").Value = 1311.7222
$ws.Range("M102").Value = 310.2778000000001
$ws.Range("H113").Value = 8073.913
$ws.Range("J113").Value = 9343.75
$ws.Range("L113").Value = 9343.75
$ws.Range("N113").Value = -13683.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3999.7778
$ws.Range("J46").Value = 3678.6428
$ws.Range("L46").Value = 3678.6428
$ws.Range("N46").Value = -4054.6428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 179077.12
$ws.Range("I122").Value = 290016.66
$ws.Range("J122").Value = 6504.5557
$ws.Range("K122").Value = 870049.98
$ws.Range("L122").Value = 19513.6671
$ws.Range("M122").Value = -867599.98
$ws.Range("N122").Value = -24413.6671
$ws.Range("H126").Value = 7999.6665
$ws.Range("J126").Value = 9000.0
$ws.Range("L126").Value = 27000.0
$ws.Range("N126").Value = -31940.0
